# ============================================================================
# Weekly CompStat report refresh: new crime data collected for the week of
# 5/20/2024 - 5/26/2024 (Volume 31, Number 21). Updates shared header text,
# refreshes all Crime Complaints figures (rows 15-31), and normalizes the
# "H" column width to match its siblings.
# ============================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report header / banner text -------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  21"
$ws.Range("C9").Value = "Report Covering the Week  5/20/2024  Through  5/26/2024"

# --- Column H width: now matches columns I/J (was previously wider) --------------
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(9).ColumnWidth

# --- Complex cells: style/type changes (copy template cell, then set value if numeric) ---
$ws.Range("C14").Copy($ws.Range("C15"))
$ws.Range("C14").Copy($ws.Range("D16"))
$ws.Range("E14").Copy($ws.Range("E16"))
$ws.Range("F15").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 1
$ws.Range("L14").Copy($ws.Range("E22"))
$ws.Range("E22").Value = -100
$ws.Range("F15").Copy($ws.Range("C23"))
$ws.Range("C23").Value = 6
$ws.Range("C14").Copy($ws.Range("C28"))
$ws.Range("F15").Copy($ws.Range("D28"))
$ws.Range("D28").Value = 1
$ws.Range("L14").Copy($ws.Range("E28"))
$ws.Range("E28").Value = -100
$ws.Range("F15").Copy($ws.Range("D31"))
$ws.Range("D31").Value = 1
$ws.Range("L14").Copy($ws.Range("E31"))
$ws.Range("E31").Value = -100

# --- Simple cells: value-only changes (style/type unchanged) ---
$ws.Range("L15").Value = 20
$ws.Range("N15").Value = -40
$ws.Range("C16").Value = 5
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 45
$ws.Range("K16").Value = -18.181818181818
$ws.Range("L16").Value = -25
$ws.Range("M16").Value = 25
$ws.Range("N16").Value = -81.781376518218
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 200
$ws.Range("F17").Value = 15
$ws.Range("H17").Value = 15.384615384615
$ws.Range("I17").Value = 50
$ws.Range("J17").Value = 56
$ws.Range("K17").Value = -10.714285714285
$ws.Range("L17").Value = 4.166666666666
$ws.Range("M17").Value = -5.660377358490
$ws.Range("N17").Value = -45.054945054945
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 5
$ws.Range("H18").Value = -50
$ws.Range("I18").Value = 35
$ws.Range("J18").Value = 47
$ws.Range("K18").Value = -25.531914893617
$ws.Range("L18").Value = -46.969696969697
$ws.Range("M18").Value = -20.454545454545
$ws.Range("N18").Value = -82.412060301507
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 166.666666666667
$ws.Range("F19").Value = 54
$ws.Range("G19").Value = 51
$ws.Range("H19").Value = 5.882352941176
$ws.Range("I19").Value = 268
$ws.Range("J19").Value = 288
$ws.Range("K19").Value = -6.944444444444
$ws.Range("L19").Value = 5.928853754940
$ws.Range("M19").Value = 15.517241379310
$ws.Range("N19").Value = -20.474777448071
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("I20").Value = 21
$ws.Range("J20").Value = 23
$ws.Range("K20").Value = -8.695652173913
$ws.Range("L20").Value = -12.5
$ws.Range("M20").Value = 16.666666666666
$ws.Range("N20").Value = -90.322580645161
$ws.Range("C21").Value = 29
$ws.Range("D21").Value = 12
$ws.Range("E21").Value = 141.666666666667
$ws.Range("F21").Value = 88
$ws.Range("G21").Value = 89
$ws.Range("H21").Value = -1.123595505617
$ws.Range("I21").Value = 425
$ws.Range("J21").Value = 471
$ws.Range("K21").Value = -9.766454352441
$ws.Range("L21").Value = -7.205240174672
$ws.Range("M21").Value = 9.819121447028
$ws.Range("N21").Value = -61.468721668177
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -50
$ws.Range("J22").Value = 17
$ws.Range("K22").Value = -23.529411764705
$ws.Range("E23").Value = 500
$ws.Range("F23").Value = 7
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 133.333333333333
$ws.Range("I23").Value = 20
$ws.Range("J23").Value = 18
$ws.Range("K23").Value = 11.111111111111
$ws.Range("L23").Value = -20
$ws.Range("M23").Value = 5.263157894736
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 11
$ws.Range("E24").Value = 109.090909090909
$ws.Range("F24").Value = 64
$ws.Range("G24").Value = 58
$ws.Range("H24").Value = 10.344827586206
$ws.Range("I24").Value = 258
$ws.Range("J24").Value = 269
$ws.Range("K24").Value = -4.089219330855
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -20.123839009287
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -33.333333333333
$ws.Range("F25").Value = 22
$ws.Range("G25").Value = 31
$ws.Range("H25").Value = -29.032258064516
$ws.Range("I25").Value = 108
$ws.Range("J25").Value = 159
$ws.Range("K25").Value = -32.075471698113
$ws.Range("L25").Value = -17.557251908396
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = -12.5
$ws.Range("F26").Value = 18
$ws.Range("G26").Value = 25
$ws.Range("H26").Value = -28
$ws.Range("I26").Value = 122
$ws.Range("J26").Value = 135
$ws.Range("K26").Value = -9.629629629629
$ws.Range("L26").Value = -5.426356589147
$ws.Range("M26").Value = -15.277777777777
$ws.Range("C27").Value = 1
$ws.Range("F27").Value = 5
$ws.Range("I27").Value = 14
$ws.Range("K27").Value = 366.666666666667
$ws.Range("L27").Value = 133.333333333333
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 28
$ws.Range("K28").Value = -17.857142857142
$ws.Range("G31").Value = 2
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 12
$ws.Range("K31").Value = -75
$ws.Range("L31").Value = -62.5
